$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert the "Website address:" block right after the "(Joanna ... Drozd)"
#    paragraph and before the existing hyperlink paragraph:
#       <empty paragraph>
#       Website address:
# ---------------------------------------------------------------------------
$introPara = $d.Paragraphs.Item(2)
$r = $introPara.Range
$r.Collapse(0)
$r.Text = [char]13
$r.Collapse(0)
$r.Text = [char]13
$r.Collapse(0)

$blankPara1 = $d.Paragraphs.Item(3)
$labelPara1 = $d.Paragraphs.Item(4)
$labelPara1.Range.Text = "Website address:"

# ---------------------------------------------------------------------------
# 2) Clean up the original hyperlink paragraph (now paragraph 5):
#    - collapse all hyperlink runs into a single run with the shortened
#      display text (drop the "/index.html" suffix from the link itself)
#    - replace the trailing "/index.html" run with a single plain space
# ---------------------------------------------------------------------------
$hyperlinkPara = $d.Paragraphs.Item(5)
$existingLink = $d.Hyperlinks.Item(1)
$existingLink.TextToDisplay = "https://moveitdanceacademy.github.io"

$paraEnd = $hyperlinkPara.Range.End
$tail = $d.Range($existingLink.Range.End, $paraEnd - 1)
if ($tail.Start -lt $tail.End) {
    $tail.Delete()
}
$insertionPoint = $d.Range($existingLink.Range.End, $existingLink.Range.End)
$insertionPoint.InsertAfter(" ")

# ---------------------------------------------------------------------------
# 3) Insert the "Repository:" block right after that paragraph:
#       <empty paragraph>
#       Repository:
#       <hyperlink to the GitHub repo>
# ---------------------------------------------------------------------------
$r2 = $hyperlinkPara.Range
$r2.Collapse(0)
$r2.Text = [char]13
$r2.Collapse(0)
$r2.Text = [char]13
$r2.Collapse(0)
$r2.Text = [char]13
$r2.Collapse(0)

$blankPara2 = $d.Paragraphs.Item(6)
$labelPara2 = $d.Paragraphs.Item(7)
$labelPara2.Range.Text = "Repository:"

$repoPara = $d.Paragraphs.Item(8)
$repoUrl = "https://github.com/moveitdanceacademy/moveitdanceacademy.github.io"
$newLink = $d.Hyperlinks.Add($repoPara.Range, $repoUrl, $null, $null, $repoUrl)

Write-Output "Final paragraph dump:"
$i = 0
foreach ($p in $d.Paragraphs) {
    $i++
    Write-Output "Para $i : [$($p.Range.Text)]"
}
